$wb = $excel.ActiveWorkbook

# --- Sheets involved -------------------------------------------------
#   Overview (sheet1) : columns E (zh-cn) and F (de-de) show the status text
#   zh-cn    (sheet2) : column C is the "Status" column
#   de-de    (sheet3) : column C is the "Status" column
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text update: "Ready for handoff" -> "In Translation" ------------
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value     = "In Translation"
$wsDeDe.Range("C2:C4").Value     = "In Translation"

# --- Column width update (status columns shrink to fit new text) -----
# NOTE: Excel quantizes ColumnWidth to a 1/6-character pixel grid before
# storing it, so the input below is chosen (12.5) so that, after that
# quantization, the persisted <col width="..."> lands on 13.3333333333333,
# the closest grid point to the target width of 13.4101845877511.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
